$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 806
$ws.Range("F4").Value = 290
$ws.Range("F5").Value = 498
$ws.Range("F6").Value = 1136
$ws.Range("F8").Value = 42
$ws.Range("F9").Value = 117
$ws.Range("F10").Value = 119
$ws.Range("F11").Value = 1168
$ws.Range("F14").Value = 811
$ws.Range("F15").Value = 835
$ws.Range("F17").Value = 56
$ws.Range("F20").Value = 700
$ws.Range("F21").Value = 1724
$ws.Range("F22").Value = 2474
$ws.Range("F23").Value = 691
$ws.Range("F24").Value = 74
$ws.Range("F25").Value = 1961
$ws.Range("F26").Value = 368
$ws.Range("F27").Value = 2844
$ws.Range("F28").Value = 527
$ws.Range("F30").Value = 695
$ws.Range("F31").Value = 133
$ws.Range("F32").Value = 107
$ws.Range("F34").Value = 987
$ws.Range("F35").Value = 1708
$ws.Range("F36").Value = 350
$ws.Range("F38").Value = 540
$ws.Range("F39").Value = 164
$ws.Range("F41").Value = 159
$ws.Range("F42").Value = 21

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F9").Value = 9
$ws.Range("F12").Value = 73

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 806
$ws.Range("F5").Value = 290
$ws.Range("F6").Value = 498
$ws.Range("F7").Value = 1136
$ws.Range("F9").Value = 42
$ws.Range("F10").Value = 117
$ws.Range("F11").Value = 119
$ws.Range("F12").Value = 1168
$ws.Range("F14").Value = 811
$ws.Range("F15").Value = 835
$ws.Range("F18").Value = 56
$ws.Range("F21").Value = 700
$ws.Range("F22").Value = 1724
$ws.Range("F23").Value = 2474
$ws.Range("F24").Value = 691
$ws.Range("F25").Value = 74
$ws.Range("F28").Value = 2844
$ws.Range("F29").Value = 527
$ws.Range("F31").Value = 9
$ws.Range("F35").Value = 73
$ws.Range("F36").Value = 695
$ws.Range("F37").Value = 133
$ws.Range("F38").Value = 107
$ws.Range("F40").Value = 987
$ws.Range("F41").Value = 1708
$ws.Range("F43").Value = 350
$ws.Range("F44").Value = 540
$ws.Range("F45").Value = 164
$ws.Range("F47").Value = 159
$ws.Range("F48").Value = 21

$wb.Save()